$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Bmp6", "Bmpr2", "ECs", 2, 0.6666666666666666, 31.825501, 95.47650300000001, 0.6063608875535647, 0.6063608875535647, 3, 1, 40.70766766666667, 122.123003, 0.3776398983502007, 0.3776398983502007, 1295.541918033168, 11659.87726229851, 0.2289860639392657, 0.2289860639392656),
    @("ECs", "Bmp6", "Bmpr2", "FAPs", 2, 0.6666666666666666, 31.825501, 95.47650300000001, 0.6063608875535647, 0.6063608875535647, 3, 1, 39.715023, 119.145069, 0.3684312589831062, 0.3684312589831062, 1263.950504201523, 11375.55453781371, 0.2234023051994735, 0.2234023051994735),
    @("ECs", "Bmp6", "Bmpr2", "sCs", 2, 0.6666666666666666, 31.825501, 95.47650300000001, 0.6063608875535647, 0.6063608875535647, 3, 1, 27.37224266666666, 82.11672799999999, 0.253928842666693, 0.253928842666693, 871.1353363602427, 7840.218027242184, 0.1539725184148255, 0.1539725184148255),
    @("FAPs", "Bmp6", "Bmpr2", "ECs", 3, 1, 4.846280666666666, 14.538842, 0.09233460445363234, 0.09233460445363234, 3, 1, 40.70766766666667, 122.123003, 0.3776398983502007, 0.3776398983502007, 197.2807827980584, 1775.527045182526, 0.03486923064007571, 0.0348692306400757),
    @("FAPs", "Bmp6", "Bmpr2", "FAPs", 3, 1, 4.846280666666666, 14.538842, 0.09233460445363234, 0.09233460445363234, 3, 1, 39.715023, 119.145069, 0.3684312589831062, 0.3684312589831062, 192.470148141122, 1732.231333270098, 0.03401895456655889, 0.03401895456655889),
    @("FAPs", "Bmp6", "Bmpr2", "sCs", 3, 1, 4.846280666666666, 14.538842, 0.09233460445363234, 0.09233460445363234, 3, 1, 27.37224266666666, 82.11672799999999, 0.253928842666693, 0.253928842666693, 132.6535704387751, 1193.882133948976, 0.02344641924699774, 0.02344641924699774),
    @("sCs", "Bmp6", "Bmpr2", "ECs", 3, 1, 15.81429, 47.44287, 0.301304507992803, 0.3013045079928031, 3, 1, 40.70766766666667, 122.123003, 0.3776398983502007, 0.3776398983502007, 643.76286170429, 5793.86575533861, 0.1137846037708594, 0.1137846037708594),
    @("sCs", "Bmp6", "Bmpr2", "FAPs", 3, 1, 15.81429, 47.44287, 0.301304507992803, 0.3013045079928031, 3, 1, 39.715023, 119.145069, 0.3684312589831062, 0.3684312589831062, 628.06489107867, 5652.584019708031, 0.1110099992170738, 0.1110099992170738),
    @("sCs", "Bmp6", "Bmpr2", "sCs", 3, 1, 15.81429, 47.44287, 0.301304507992803, 0.3013045079928031, 3, 1, 27.37224266666666, 82.11672799999999, 0.253928842666693, 0.253928842666693, 432.87258348104, 3895.85325132936, 0.07650990500486983, 0.07650990500486984)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
    }
}
